$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Challenges")

# New header cells for the Gameweeks import feature.
$ws.Range("S1").Value = "Show Statistics Continuously"
$ws.Range("T1").Value = "Gameweek"

# New data row values.
# "true" must land as literal text (not Excel's auto-detected boolean), so
# it is entered as a formula producing the string "true" and then converted
# to a plain value via copy / paste-special (values only). This avoids both
# the boolean auto-coercion that a direct Value="true" assignment triggers
# and the quote-prefix cell style that a leading-apostrophe text entry
# would otherwise leave behind.
$s2 = $ws.Range("S2")
$s2.Formula = "=""true"""
$s2.Copy()
$s2.PasteSpecial(-4163)  # xlPasteValues

$ws.Range("T2").Value = 1
